$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.523.13'
$ws.Range('E2').Value = '  +4.06%  '
$ws.Range('D3').Value = '1.737.03'
$ws.Range('E3').Value = '  +4.47%  '
$ws.Range('E4').Value = '  +0.11%  '
$c = $ws.Range('D5')
$c.Value = "'244.22"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +3.74%  '
$ws.Range('E6').Value = '  +0.08%  '
$c = $ws.Range('D7')
$c.Value = "'0.4798"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +3.53%  '
$c = $ws.Range('D8')
$c.Value = "'0.2668"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.63%  '
$c = $ws.Range('D9')
$c.Value = "'0.06229"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').Value = '1.736.74'
$ws.Range('E10').Value = '  +4.70%  '
$ws.Range('E11').Value = '  +2.66%  '
$c = $ws.Range('D12')
$c.Value = "'15.75"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +7.68%  '
$c = $ws.Range('D13')
$c.Value = "'0.6169"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +7.79%  '
$c = $ws.Range('D14')
$c.Value = "'4.533"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +4.74%  '
$c = $ws.Range('D15')
$c.Value = "'76.83"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +2.45%  '
$c = $ws.Range('D16')
$c.Value = "'1.001"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '26.521.57'
$ws.Range('E17').Value = '  +4.06%  '
$ws.Range('E18').Value = '  +0.07%  '
$c = $ws.Range('D19')
$c.Value = "'0.000006898"
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.Value = "'11.73"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +3.33%  '
$ws.Range('D21').Value = '1.960.47'
$ws.Range('E21').Value = '  +4.49%  '
$c = $ws.Range('D22')
$c.Value = "'4.562"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +3.41%  '
$c = $ws.Range('D23')
$c.Value = "'8.890"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +3.12%  '
$c = $ws.Range('D24')
$c.Value = "'5.329"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.14%  '
$c = $ws.Range('D25')
$c.Value = "'135.66"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.39%  '
$c = $ws.Range('D26')
$c.Value = "'15.35"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +3.33%  '
$c = $ws.Range('D27')
$c.Value = "'1.799"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +5.31%  '
$c = $ws.Range('D28')
$c.Value = "'1.409"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +3.14%  '
$c = $ws.Range('D29')
$c.Value = "'106.56"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +2.89%  '
$c = $ws.Range('D30')
$c.Value = "'3.987"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.87%  '
$c = $ws.Range('D31')
$c.Value = "'3.714"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +3.29%  '
$c = $ws.Range('D32')
$c.Value = "'0.07880"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +2.17%  '
$c = $ws.Range('D33')
$c.Value = "'0.04565"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +5.42%  '
$c = $ws.Range('D34')
$c.Value = "'2.616"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D35')
$c.Value = "'0.6348"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +5.79%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D36')
$c.Value = "'0.9946"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +5.78%  '
$c = $ws.Range('D37')
$c.Value = "'0.9315"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +2.23%  '
$c = $ws.Range('D38')
$c.Value = "'111.05"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.72%  '
$c = $ws.Range('D39')
$c.Value = "'2.445"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.07%  '
$c = $ws.Range('D40')
$c.Value = "'1.980"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +8.40%  '
$c = $ws.Range('D41')
$c.Value = "'1.005"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.63%  '
$c = $ws.Range('D42')
$c.Value = "'0.01511"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.92%  '
$c = $ws.Range('D43')
$c.Value = "'5.691"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +13.98%  '
$c = $ws.Range('D44')
$c.Value = "'0.3901"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +5.14%  '
$c = $ws.Range('D45')
$c.Value = "'6.893"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +12.71%  '
$c = $ws.Range('D46')
$c.Value = "'0.1193"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +7.88%  '
$c = $ws.Range('D47')
$c.Value = "'0.05333"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.51%  '
$c = $ws.Range('D48')
$c.Value = "'7.907"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.92%  '
$c = $ws.Range('D49')
$c.Value = "'30.81"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +2.15%  '
$c = $ws.Range('D50')
$c.Value = "'1.255"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +5.44%  '
$c = $ws.Range('D51')
$c.Value = "'0.3435"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.92%  '
